# pool-importvorlage-komplett.xlsx — "add new fields and corrections"
#
# Adds three new import columns to the product pool template:
#   - isEcoProduct            (inserted right after "brand", before "categories1")
#   - unit1 / minimumQuantity1 / maximumQuantity1
#                              (inserted right after "inventory", before "specialPrice")
# and fills in the matching sample values in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column: isEcoProduct -------------------------------------------------
# Currently: ... AA=manufacturer AB=brand AC=categories1 ...
# Insert a blank column at AC so categories1..specialPrice shift one to the right.
$ws.Range("AC1").EntireColumn.Insert() | Out-Null
$ws.Range("AC1").Value = "isEcoProduct"
$ws.Range("AC2").Value = "0 oder leer = Nein, 1 = Ja"

# --- New columns: unit1, minimumQuantity1, maximumQuantity1 ------------------
# After the previous insert: ... AU=stockStatus AV=inventory AW=specialPrice ...
# Insert three blank columns at AW so specialPrice/startDate/endDate shift right by 3.
$ws.Range("AW1:AY1").EntireColumn.Insert() | Out-Null

$ws.Range("AW1").Value = "unit1"
$ws.Range("AX1").Value = "minimumQuantity1"
$ws.Range("AY1").Value = "maximumQuantity1"

$ws.Range("AW2").Value = "Verpackungseinheit1 "
$ws.Range("AX2").Value = "mind. Anzahl Bestelleinheiten"
$ws.Range("AY2").Value = "max. Anzahl Bestelleinheiten"

# Leave the selection where the author ended up editing.
$ws.Range("AY2").Select() | Out-Null
